$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEnglishStory = @"
Winter arrived quietly—too quietly. One morning, the town woke up under a pale sky, and everyone immediately complained. The air was sharp, clean, and rude, like it had no manners at all. Frost clung to windows, drawing fancy ice patterns that looked impressive but were absolutely useless. Every breath turned into a little cloud, as if people had suddenly become part-time dragons.
The trees stood bare and silent, probably because they had given up arguing with winter. Their branches reached upward like they were asking, “Excuse me, when is spring coming back?” Roads glittered with frost, and every step made a crunchy sound that reminded people they had chosen the wrong shoes. Children rushed outside in oversized coats, moving like tiny walking pillows, laughing as they threw snowballs that mostly missed their targets.
Inside the houses, winter acted very differently. Fires crackled happily, filling rooms with warm light and the smell of tea, spices, and slightly burnt toast. Windows fogged up instantly, turning the outside world into a blurry mystery. Time slowed down, not because winter was peaceful, but because everyone moved slower while wearing five layers of clothing.
At night, snow began to fall—soft, steady, and silent, like it was sneaking in uninvited. Under streetlights, snowflakes sparkled dramatically before disappearing, clearly showing off. The streets looked magical, calm, and slippery enough to cause at least one embarrassing fall.
Winter, it turned out, wasn’t just about cold—it was about survival. Beneath the frozen ground, life waited patiently, probably wearing socks. Though the land slept, hope stayed awake, counting the days until spring. Until then, winter hugged the world tightly, refusing to let go, reminding everyone that even in the coldest moments… hot tea and thick blankets are true heroes.
Inside the houses, winter felt warmer. Fires crackled softly, filling rooms with golden light. The scent of tea and spices lingered in the air, and windows fogged as stories were shared. Time slowed during winter, as if the season itself encouraged people to pause. Evenings stretched longer, inviting quiet moments, thoughtful conversations, and peaceful reflection.
At night, snow began to fall—slow, steady, and silent. Under streetlights, each snowflake shimmered briefly before disappearing, turning ordinary streets into something magical. The world seemed calmer, wrapped in a gentle stillness.
Winter was not only a season of cold, but of rest. Beneath frozen soil, life waited patiently. Though the land slept, hope remained alive, knowing that spring would return. Until then, winter held the world in a quiet embrace, reminding everyone that beauty can exist even in the coldest moments.
"@

$newBanglaStory = @"
শীত এসেছিল খুব চুপিচুপি—এতটাই চুপিচুপি যে সন্দেহ হওয়াই স্বাভাবিক। এক সকালে শহরটি ফ্যাকাসে আকাশের নিচে জেগে উঠল, আর সবাই একসাথে অভিযোগ করা শুরু করল। বাতাস ছিল তীক্ষ্ণ, পরিষ্কার এবং বেশ বেয়াদব—যেন তার কোনো ভদ্রতা নেই। জানালার কাঁচে তুষার লেগে সুন্দর সুন্দর নকশা আঁকছিল, দেখতে দারুণ হলেও কাজে একদমই আসত না। প্রতিটি নিঃশ্বাস ছোট ছোট মেঘে পরিণত হচ্ছিল, যেন মানুষ হঠাৎ করেই খণ্ডকালীন ড্রাগন হয়ে গেছে।
গাছগুলো দাঁড়িয়ে ছিল নগ্ন ও নীরব, সম্ভবত শীতের সাথে তর্ক করে হাল ছেড়ে দিয়েছে। তাদের ডালপালা আকাশের দিকে বাড়ানো, যেন বলছে, “মাফ করবেন, বসন্তটা কবে ফিরবে?” রাস্তাগুলো তুষারে ঝিলমিল করছিল, আর প্রতিটি পদক্ষেপে এমন খচখচ শব্দ হচ্ছিল যা মানুষকে মনে করিয়ে দিচ্ছিল—ভুল জুতা পরা হয়েছে। বাচ্চারা মোটা মোটা কোট পরে বাইরে ছুটে গেল, দেখতে ছোট ছোট হাঁটা বালিশের মতো, হাসতে হাসতে এমন সব স্নোবল ছুঁড়ছিল যা বেশিরভাগ সময় লক্ষ্যভ্রষ্ট হচ্ছিল।
ঘরের ভেতরে শীতের আচরণ ছিল একদম আলাদা। চুলোর আগুন খুশিমনে জ্বলছিল, ঘর ভরে উঠছিল উষ্ণ আলো আর চা, মশলা ও সামান্য পুড়ে যাওয়া টোস্টের গন্ধে। জানালাগুলো সাথে সাথে কুয়াশাচ্ছন্ন হয়ে গেল, বাইরের দুনিয়াকে ঝাপসা রহস্যে পরিণত করল। সময় ধীরে চলছিল—শীত শান্ত বলে নয়, বরং সবাই পাঁচ স্তরের কাপড় পরে ধীরে নড়াচড়া করছিল বলে।
রাতে তুষার পড়তে শুরু করল—নরম, স্থির আর নিঃশব্দে, যেন দাওয়াত ছাড়াই ঢুকে পড়ছে। রাস্তার বাতির নিচে প্রতিটি তুষারকণা নাটকীয়ভাবে ঝিলমিল করছিল, তারপর মিলিয়ে যাচ্ছিল—একেবারে নজর কাড়ার চেষ্টা। রাস্তাগুলো দেখতে জাদুকরী, শান্ত, আর যথেষ্ট পিচ্ছিল—যাতে অন্তত একজন লজ্জাজনকভাবে পিছলে পড়ে।
শেষ পর্যন্ত বোঝা গেল, শীত শুধু ঠান্ডার ঋতু নয়—এটা টিকে থাকার সময়। জমে থাকা মাটির নিচে জীবন ধৈর্য ধরে অপেক্ষা করছিল, সম্ভবত মোজা পরে। ভূমি ঘুমিয়ে থাকলেও আশা জেগে ছিল, বসন্ত আসার দিন গুনছিল। ততদিন পর্যন্ত শীত পৃথিবীটাকে শক্ত করে জড়িয়ে ধরে রাখল, ছাড়তে চাইছিল না, আর সবাইকে মনে করিয়ে দিল—সবচেয়ে ঠান্ডা মুহূর্তেও… গরম চা আর মোটা কম্বলই আসল নায়ক।
"@

$ws.Range("D3").Value = $newEnglishStory
$ws.Range("E3").Value = $newBanglaStory

$ws.Rows.Item(3).RowHeight = 246

$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("E3").Select() | Out-Null
